# Apply PivotFilterOverrides-related ToDo list updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 25 ("Issue 3" task) is now complete -> style changes from Neutral to Good
$ws.Range("B25:E25").Style = "Good"

# 2. Remove the "Issue 1: calculationType ..." row (row 26) entirely - it shifts everything below up by one
$ws.Rows.Item(26).Delete()

# After the delete:
#   row 27 -> "check export to Excel..." (formerly row 28)
#   row 28 -> "document the changes..." (formerly row 29)
#   row 29 -> "review the content..." (formerly row 30)
#   row 30 -> "(1) Data type support..." (formerly row 31)
#   row 31 -> "In NEWS.md, move the examples..." (formerly row 32)
#   row 32 -> "Add note to the docs..." (formerly row 33)

# 3. Structurally insert a new (blank, Neutral-styled) row before "In NEWS.md..." (row 31)
#    for the new filter-overrides note
$ws.Rows.Item(31).Insert()

# After this insert:
#   row 31 -> blank (new)
#   row 32 -> "In NEWS.md, move the examples..."
#   row 33 -> "Add note to the docs..."

# 4. Structurally append two new rows at the end of the list (after "Add note to the docs...", row 33)
$ws.Rows.Item(34).Insert()
$ws.Rows.Item(35).Insert()
$ws.Range("B33:F33").Copy()
$ws.Range("B34:F35").PasteSpecial(-4122)

# 5. Now fill in the cell values/text, in the same order the original author typed them
#    (this controls the order new strings are appended to the shared-strings table)

# Extend the "Add quick pivot tests" note (row 23) and add two new formatted (blank) cells F23:G23
$ws.Range("B23").Value = "Add quick pivot tests (see trello), add tests for specifying style info (baseStyleName, styleDeclarations) upfront, add tests for filteroverrides"
$ws.Range("F23:G23").Style = "Neutral"

# New rows at the bottom of the list
$ws.Range("B34").Value = "Follow normal dev workflow."
$ws.Range("B35").Value = "Release as v1.0 to CRAN."

# New row inserted above "In NEWS.md..."
$ws.Range("B31").Value = "(5) Upfront styling, (6) overriding filters as part of calculations."

# Add a trailing comma to the "(1) Data type support..." note (now at row 30)
$ws.Range("B30").Value = "(1) Data type support, (2) data type formatting, (3) Dealing with columns with illegal names, (4) New export options,"

# 6. Restore the active-cell selection to match the edited state
$ws.Range("B29").Select()
